# Single-column results table: one stat per row, each cell holding a
# single run (some rows later in the doc previously packed a whole
# tab-separated summary line into one run/cell).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $rowIndex, $text) {
    $cell = $table.Cell($rowIndex, 1)
    $r = $cell.Range
    # Exclude the end-of-cell marker so we only replace the cell's content.
    $r.End = $r.End - 1
    $r.Text = $text
}

# Rows 1-4: headline percentages / memory figures
Set-CellText $t 1  "0M"
Set-CellText $t 2  "0M"
Set-CellText $t 3  "0M"
Set-CellText $t 4  "80"

# Rows 6-12: per-phase timing stats, recomputed
Set-CellText $t 6  "0.00091"
Set-CellText $t 7  "0.00027"
Set-CellText $t 8  "0.00009"
Set-CellText $t 9  "0.00036"
Set-CellText $t 10 "0.00047"
Set-CellText $t 11 "0.00053"
Set-CellText $t 12 "0.02155"

# Rows 44-46: previously multi-run tab-separated summary lines, now
# collapsed down to the single headline value each (matching rows 1-3).
Set-CellText $t 44 "99.97"
Set-CellText $t 45 "0.02"
Set-CellText $t 46 "70"
